$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the first strategy test row with the new date/TP/SL values.
$ws.Range("D2").Value2 = 44531
$ws.Range("G2").Value2 = 0.8
$ws.Range("H2").Value2 = 0.8

# Remove the second test row's details (keep the empty K3 cell/style in place).
$ws.Range("A3:J3").Clear()

# Remove the leftover blank placeholder rows below the test data.
$ws.Range("B4:I17").Clear()

# Match the final cell selection left active on the sheet.
$null = $ws.Range("B9").Select()
